# Deploy the implementation guide.
#
# Changes applied to docs/ValueSet-age-at-onset-vs.xlsx:
#  1. On the "Metadata" sheet:
#     - "Date" value (B8) updated to a newer generation timestamp.
#     - "Contact" value (B10) updated to the resolved publisher contact.
#     - A new "Jurisdiction" property row is inserted right after "Contact",
#       pushing "Description" / "Purpose" / "Copyright" / "Immutable" down
#       by one row.
#  2. The second sheet (the "Include ..." CodeSystem sheet) is renamed to
#     "Include #0"; its data is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Make room for the new "Jurisdiction" row right after "Contact" (row 10) ---
# Copy the formatting of the last existing row down into the new row 15 so the
# appended row keeps the same style as the rest of the table body.
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)

# Shift the "Description"/"Purpose"/"Copyright"/"Immutable" rows down by one,
# working from the bottom up so nothing gets clobbered before it is read.
$ws.Range("A15").Value2 = $ws.Range("A14").Value2
$ws.Range("B15").Value2 = $ws.Range("B14").Value2

$ws.Range("A14").Value2 = $ws.Range("A13").Value2
$ws.Range("B14").Value2 = $ws.Range("B13").Value2

$ws.Range("A13").Value2 = $ws.Range("A12").Value2
$ws.Range("B13").Value2 = $ws.Range("B12").Value2

$ws.Range("A12").Value2 = $ws.Range("A11").Value2
$ws.Range("B12").Value2 = $ws.Range("B11").Value2

# --- Insert the new "Jurisdiction" property (no value) in the freed row 11 ---
$ws.Range("A11").Value2 = "Jurisdiction"
$ws.Range("B11").Value2 = ""

# --- Refresh the generation timestamp and contact display string ---
$ws.Range("B8").Value2 = "2024-10-02T15:04:17+00:00"
$ws.Range("B10").Value2 = "Ferlab.bio (http://example.org/example-publisher)"

# --- Rename the CodeSystem "Include" sheet tab ---
$ws2 = $wb.Worksheets.Item("Include from Ferlab.bio CodeS")
$ws2.Name = "Include #0"
